$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

# Map of row number -> Result value, matching the rows that were updated
# in the DSL / screenshot rework described in the commit message.
$results = [ordered]@{
    2  = "Pass"
    3  = "Pass"
    4  = "Pass"
    5  = "Pass"
    6  = "Pass"
    7  = "Pass"
    11 = "Pass"
    12 = "Pass"
    13 = "Pass"
    14 = "Pass"
    15 = "Fail"
    16 = "Pass"
    17 = "Pass"
    18 = "Fail"
    19 = "Fail"
    20 = "Fail"
    21 = "Fail"
    22 = "Pass"
    23 = "Pass"
    24 = "Pass"
    25 = "Pass"
}

foreach ($row in $results.Keys) {
    $ws.Range("J$row").Value = $results[$row]
}
